# Update the "Metadata" worksheet (sheet 1) of the ValueSet workbook:
#  - bump Version 0.1.6 -> 0.1.7
#  - Status active -> draft
#  - Date updated
#  - Contact rows updated (publisher contact + Bob Milius)
#  - a new "Jurisdiction" row is inserted before "Description", pushing
#    Description/Purpose/Copyright/Immutable down by one row each

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- simple in-place metadata updates -------------------------------------
$ws.Range("B3").Value = "0.1.7"
$ws.Range("B6").Value = "draft"
$ws.Range("B8").Value = "2024-11-22T12:33:30-06:00"
$ws.Range("B10").Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"
$ws.Range("B11").Value = "Bob Milius (bmilius@nmdp.org)"

# --- capture the existing Description/Purpose/Copyright/Immutable block ---
# before we start overwriting it, so we can shift it down by one row.
$a12 = $ws.Range("A12").Value()
$b12 = $ws.Range("B12").Value()
$a13 = $ws.Range("A13").Value()
$b13 = $ws.Range("B13").Value()
$a14 = $ws.Range("A14").Value()
$b14 = $ws.Range("B14").Value()
$a15 = $ws.Range("A15").Value()
$b15 = $ws.Range("B15").Value()

# --- shift rows 12-15 down to 13-16 (bottom-up so nothing is clobbered) ---
$ws.Range("A16").Value = $a15
$ws.Range("B16").Value = $b15

$ws.Range("A15").Value = $a14
$ws.Range("B15").Value = $b14

$ws.Range("A14").Value = $a13
$ws.Range("B14").Value = $b13

$ws.Range("A13").Value = $a12
$ws.Range("B13").Value = $b12

# --- insert the new Jurisdiction row at row 12 -----------------------------
$ws.Range("A12").Value = "Jurisdiction"
$ws.Range("B12").Value = ""
